$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41: add D41 (date) and E41 (text) ---
$ws.Range("B41").Copy()
$ws.Range("D41").PasteSpecial(-4122)   # xlPasteFormats - reuse the date number format/style
$ws.Range("D41").Value = 43724
$ws.Range("E41").Value = "Now it is meta_table_example"

# --- Row 42: add D42 (date) and E42 (text) ---
$ws.Range("B41").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("D42").Value = 43724
$ws.Range("E42").Value = "added a get_code_tree() function and improved the help page"

# --- Row 43: new row with A43, B43 (date), C43 ---
$ws.Range("C43").Value = "Option in write_ggplot to include survival plot from ggsurvplot - package survminer"
$ws.Range("A43").Value = "JK"
$ws.Range("B41").Copy()
$ws.Range("B43").PasteSpecial(-4122)
$ws.Range("B43").Value = 43713

$excel.CutCopyMode = 0

# --- Update view: scroll position and active selection ---
$ws.Application.GoTo($ws.Range("A31"))
$ws.Range("D43").Select()
